# Update the "2024" worksheet: a new September log entry was added at the
# top of the combined "details / date" log list that runs down column
# R/S (September_Details / September_Date, rows 35-97) and then continues
# into column P/Q (August_Details / August_Date, rows 98-101). Adding one
# entry at the top pushes every later entry down by one row; since the
# R/S portion of the list can only hold 63 rows (35-97) before the sheet
# layout hands off to the P/Q portion, the list boundary itself moves down
# one row too (R/S now fills rows 35-98, P/Q now fills rows 99-102).
# Finally the "Broadband" label that used to sit right after the list in
# column A (row 102) is pushed down to row 103.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$firstRow = 35          # first row of the R/S list
$rsLastRowBefore = 97   # last row of the R/S list, before the edit
$pqLastRowBefore = 101  # last row of the P/Q list, before the edit
$labelRowBefore = 102   # row holding the "Broadband" label, before the edit

# --- Collect current (before-edit) values of the combined list ---
$values = @()
$values += , @("corporate internet share", "2024-09-09 11:17:34")

for ($r = $firstRow; $r -le $rsLastRowBefore; $r++) {
    $detail = $ws.Range("R$r").Value2
    $date = $ws.Range("S$r").Value2
    $values += , @($detail, $date)
}

for ($r = ($rsLastRowBefore + 1); $r -le $pqLastRowBefore; $r++) {
    $detail = $ws.Range("P$r").Value2
    $date = $ws.Range("Q$r").Value2
    $values += , @($detail, $date)
}

# The label that used to sit right after the list, in column A
$label = $ws.Range("A$labelRowBefore").Value2

$rsCountBefore = $rsLastRowBefore - $firstRow + 1   # 63
$pqCountBefore = $pqLastRowBefore - $rsLastRowBefore # 4

$rsLastRowAfter = $rsLastRowBefore + 1   # 98
$pqLastRowAfter = $pqLastRowBefore + 1   # 102
$labelRowAfter = $labelRowBefore + 1     # 103

# --- Write the (now one-longer) list back out, shifted down by one row ---
$idx = 0
for ($r = $firstRow; $r -le $rsLastRowAfter; $r++) {
    $pair = $values[$idx]
    $ws.Range("R$r").Value = $pair[0]
    $ws.Range("S$r").Value = $pair[1]
    # This row used to (possibly) be a P/Q row; make sure it no longer is.
    $ws.Range("P$r").Value = ""
    $ws.Range("Q$r").Value = ""
    $idx++
}

for ($r = ($rsLastRowAfter + 1); $r -le $pqLastRowAfter; $r++) {
    $pair = $values[$idx]
    $ws.Range("P$r").Value = $pair[0]
    $ws.Range("Q$r").Value = $pair[1]
    $idx++
}

# Clear the old label row's column A and move the label one row down
$ws.Range("A$labelRowBefore").Value = ""
$ws.Range("A$labelRowAfter").Value = $label
